$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = 2
$ws.Range("D2").Value = 4
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 6
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 7
$ws.Range("I2").Value = 2
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 30
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("R2").Value = 1

# --- Row 4 updates ---
$ws.Range("C4").Value = 3
$ws.Range("E4").Value = 10
$ws.Range("P4").Value = 1

# --- Row 5 updates ---
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 5
$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 6
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = "7(could be morhen)"
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 7
$ws.Range("M5").Value = 7
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 8
$ws.Range("P5").Value = 8
$ws.Range("Q5").Value = 8
$ws.Range("R5").Value = 5
$ws.Range("S5").Value = 6
$ws.Range("T5").Value = 7
$ws.Range("U5").Value = 8

# --- View state: select U5 (matches final activeCell/sqref in the diff) ---
$ws.Range("U5").Select() | Out-Null
